$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B20").Value = 43683
$ws.Range("B20").NumberFormat = "m/d/yyyy"
$ws.Range("C20").Value = 0.5
$ws.Range("C20").NumberFormat = "h:mm"
